$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6125495433807373
$ws.Range("B1").Value = 1.081016182899475
$ws.Range("C1").Value = 4.315365791320801
$ws.Range("D1").Value = 4.67966890335083
$ws.Range("E1").Value = 1.922793507575989
